# Apply the metadata restructuring described by the diff:
# - Row 1 human-readable headers are replaced with new labels/order.
# - Row 2 (sdmx-dimension / iaest-measure annotations) reshuffled.
# - Row 3 ("dim"/"medida" annotations) reshuffled.
# - Row 4 (datatype annotations) reshuffled, with xsd:int split into
#   xsd:double (Personas) and xsd:string (Nivel estudios).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Nivel estudios, código"
$ws.Range("B1").Value = "Personas"
$ws.Range("C1").Value = "Municipio codigo"
$ws.Range("D1").Value = "Nivel estudios"
$ws.Range("E1").Value = "Municipio nombre"

$ws.Range("A2").Value = "null"
$ws.Range("B2").Value = "iaest-measure:personas"
$ws.Range("C2").Value = "null"
$ws.Range("D2").Value = "iaest-measure:nivel-estudios"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

$ws.Range("A3").Value = "null"
$ws.Range("B3").Value = "medida"
$ws.Range("C3").Value = "null"
$ws.Range("D3").Value = "medida"
$ws.Range("E3").Value = "dim"

$ws.Range("A4").Value = "null"
$ws.Range("B4").Value = "xsd:double"
$ws.Range("C4").Value = "null"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "URI-Municipio"
